# Atualização automática do relatório BI
#
# The report has one "Tempo total na fase X (dias)" column per pipeline
# phase (Backlog -> AK, Construção no Canvas -> AN, Validação -> AQ,
# Publicar na plataforma -> AT, Concluído -> AW). For every ticket (row)
# that column still in progress (i.e. the one matching the ticket's
# current phase in column "Fase atual"/C) is a live "days elapsed"
# figure. Re-running the automation simply advances the reference "now"
# a little, so every one of those still-open duration cells grows by the
# same fixed amount.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Maps the literal text in column C ("Fase atual") to the worksheet
# column that holds that phase's running duration.
$phaseColumn = @{
    "Backlog"                      = "AK"
    "Construção no Canvas"         = "AN"
    "Validação"                    = "AQ"
    "Publicar na plataforma"       = "AT"
    "Concluído"                    = "AW"
}

# Elapsed time (in days) between the previous and the new report run.
$elapsed = 0.041794

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $phase = $ws.Range("C$row").Value2
    if ([string]::IsNullOrEmpty($phase)) {
        continue
    }

    $col = $phaseColumn[$phase]
    if (-not $col) {
        continue
    }

    $cell = $ws.Range("$col$row")
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = $current + $elapsed
    }
}
